$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("Citywide Totals")
$ws.Range("K2").Value = 5641
$ws.Range("K3").Value = 5775
$ws.Range("K4").Value = 1201
$ws.Range("K5").Value = 412
$ws.Range("K6").Value = 6404
$ws.Range("K7").Value = 19433

$ws = $wb.Sheets.Item("Norwood Park")
$ws.Range("K6").Value = 10
$ws.Range("K7").Value = 44

$ws = $wb.Sheets.Item("Logan Square")
$ws.Range("K2").Value = 62
$ws.Range("K7").Value = 242

$ws = $wb.Sheets.Item("Austin")
$ws.Range("K3").Value = 392
$ws.Range("K7").Value = 1285

$ws = $wb.Sheets.Item("South Chicago")
$ws.Range("K2").Value = 146
$ws.Range("K7").Value = 432

$ws = $wb.Sheets.Item("Garfield Park")
$ws.Range("K2").Value = 226
$ws.Range("K3").Value = 305
$ws.Range("K4").Value = 40
$ws.Range("K5").Value = 20
$ws.Range("K6").Value = 245
$ws.Range("K7").Value = 836

$ws = $wb.Sheets.Item("West Pullman")
$ws.Range("K6").Value = 76
$ws.Range("K7").Value = 328

$ws = $wb.Sheets.Item("Grand Crossing")
$ws.Range("K2").Value = 190
$ws.Range("K3").Value = 217
$ws.Range("K6").Value = 191
$ws.Range("K7").Value = 659

$ws = $wb.Sheets.Item("New City")
$ws.Range("K2").Value = 146
$ws.Range("K7").Value = 447

$ws = $wb.Sheets.Item("By Neighborhood")
$ws.Range("K5").Value = 48
$ws.Range("K7").Value = 570
$ws.Range("K8").Value = 1285
$ws.Range("K10").Value = 107
$ws.Range("K11").Value = 369
$ws.Range("K15").Value = 200
$ws.Range("K20").Value = 456
$ws.Range("K22").Value = 54
$ws.Range("K23").Value = 200
$ws.Range("K25").Value = 91
$ws.Range("K29").Value = 1052
$ws.Range("K31").Value = 211
$ws.Range("K33").Value = 836
$ws.Range("K34").Value = 107
$ws.Range("K37").Value = 659
$ws.Range("K42").Value = 720
$ws.Range("K43").Value = 167
$ws.Range("K46").Value = 41
$ws.Range("K48").Value = 248
$ws.Range("K50").Value = 94
$ws.Range("K53").Value = 242
$ws.Range("K54").Value = 373
$ws.Range("K55").Value = 214
$ws.Range("K63").Value = 58
$ws.Range("K65").Value = 447
$ws.Range("K67").Value = 750
$ws.Range("K69").Value = 44
$ws.Range("K73").Value = 172
$ws.Range("K76").Value = 266
$ws.Range("K79").Value = 489
$ws.Range("K83").Value = 432
$ws.Range("K84").Value = 151
$ws.Range("K85").Value = 915
$ws.Range("K88").Value = 211
$ws.Range("K89").Value = 286
$ws.Range("K90").Value = 177
$ws.Range("K93").Value = 71
$ws.Range("K94").Value = 262
$ws.Range("K95").Value = 328
$ws.Range("K101").Value = 19433

$ws = $wb.Sheets.Item("Gage Park")
$ws.Range("K3").Value = 54
$ws.Range("K7").Value = 211

$ws = $wb.Sheets.Item("North Lawndale")
$ws.Range("K2").Value = 215
$ws.Range("K3").Value = 269
$ws.Range("K6").Value = 209
$ws.Range("K7").Value = 750

$ws = $wb.Sheets.Item("South Deering")
$ws.Range("K3").Value = 57
$ws.Range("K6").Value = 30
$ws.Range("K7").Value = 151

$ws = $wb.Sheets.Item("Loop")
$ws.Range("K2").Value = 61
$ws.Range("K6").Value = 199
$ws.Range("K7").Value = 373

$ws = $wb.Sheets.Item("Englewood")
$ws.Range("K2").Value = 300
$ws.Range("K3").Value = 375
$ws.Range("K6").Value = 302
$ws.Range("K7").Value = 1052

$ws = $wb.Sheets.Item("Lake View")
$ws.Range("K3").Value = 59
$ws.Range("K7").Value = 248

$ws = $wb.Sheets.Item("River North")
$ws.Range("K2").Value = 58
$ws.Range("K6").Value = 139
$ws.Range("K7").Value = 266

$ws = $wb.Sheets.Item("Humboldt Park")
$ws.Range("K2").Value = 193
$ws.Range("K7").Value = 720

$ws = $wb.Sheets.Item("Avondale")
$ws.Range("K2").Value = 31
$ws.Range("K4").Value = 7
$ws.Range("K7").Value = 107

$ws = $wb.Sheets.Item("Lower West Side")
$ws.Range("K2").Value = 66
$ws.Range("K7").Value = 214

$ws = $wb.Sheets.Item("Jefferson Park")
$ws.Range("K2").Value = 17
$ws.Range("K7").Value = 41

$ws = $wb.Sheets.Item("Douglas")
$ws.Range("K3").Value = 70
$ws.Range("K7").Value = 200

$ws = $wb.Sheets.Item("Roseland")
$ws.Range("K3").Value = 157
$ws.Range("K7").Value = 489

$ws = $wb.Sheets.Item("Near South Side")
$ws.Range("K2").Value = 26
$ws.Range("K6").Value = 46

$ws = $wb.Sheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 149
$ws.Range("K3").Value = 149
$ws.Range("K7").Value = 456

$ws = $wb.Sheets.Item("West Lawn")
$ws.Range("K3").Value = 16
$ws.Range("K7").Value = 71

$ws = $wb.Sheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 193
$ws.Range("K7").Value = 570

$ws = $wb.Sheets.Item("Garfield Ridge")
$ws.Range("K2").Value = 39
$ws.Range("K7").Value = 107

$ws = $wb.Sheets.Item("West Loop")
$ws.Range("K3").Value = 53
$ws.Range("K7").Value = 262

$ws = $wb.Sheets.Item("East Side")
$ws.Range("K2").Value = 31
$ws.Range("K7").Value = 91

$ws = $wb.Sheets.Item("Brighton Park")
$ws.Range("K2").Value = 72
$ws.Range("K7").Value = 200

$ws = $wb.Sheets.Item("Lincoln Square")
$ws.Range("K2").Value = 25
$ws.Range("K7").Value = 94

$ws = $wb.Sheets.Item("Belmont Cragin")
$ws.Range("K2").Value = 128
$ws.Range("K7").Value = 369

$ws = $wb.Sheets.Item("Portage Park")
$ws.Range("K6").Value = 60
$ws.Range("K7").Value = 172

$ws = $wb.Sheets.Item("United Center")
$ws.Range("K6").Value = 88
$ws.Range("K7").Value = 211

$ws = $wb.Sheets.Item("Uptown")
$ws.Range("K2").Value = 80
$ws.Range("K3").Value = 89
$ws.Range("K6").Value = 85
$ws.Range("K7").Value = 286

$ws = $wb.Sheets.Item("Armour Square")
$ws.Range("K6").Value = 23
$ws.Range("K7").Value = 48

$ws = $wb.Sheets.Item("Washington Heights")
$ws.Range("K2").Value = 63
$ws.Range("K7").Value = 177

$ws = $wb.Sheets.Item("Hyde Park")
$ws.Range("K4").Value = 22
$ws.Range("K7").Value = 167

$ws = $wb.Sheets.Item("South Shore")
$ws.Range("K2").Value = 300
$ws.Range("K3").Value = 309
$ws.Range("K7").Value = 915

$ws = $wb.Sheets.Item("Clearing")
$ws.Range("K2").Value = 25
$ws.Range("K7").Value = 54
